$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Created the token generator"
$ws.Range("B9").Value = 3

$ws.Range("A10").Value = "Created the images for the token shapes and numbers"
$ws.Range("B10").Value = 3

$ws.Range("A11").Value = "Logical Architecture Diagram"
$ws.Range("B11").Value = 3

$ws.Range("B9:B11").HorizontalAlignment = -4108

$ws.Range("B12").Select() | Out-Null
